# Refreshed cohort definitions from PL 3.32.0. Caused some cohortIDs to change.
# The "targets" sheet's cohort_definition_id for the
# "Pulmonary hypertension associated with lung diseases and/or hypoxia (WHO Group 3)"
# row changed from 751 to 1265.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("targets")

$ws.Range("A4").Value = 1265
